$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'jeca'
$ws.Range("B2").Value = 'milev'
$ws.Range("C2").Value = 'kdkdljfs'
$ws.Range("D2").Value = 'lkfdsf'
$ws.Range("E2").Formula = '=TEXT(135,"0")'
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4163) | Out-Null

$ws.Range("A3").Value = 'ghsd'
$ws.Range("B3").Value = 'kjsdksjd'
$ws.Range("C3").Value = 'snjksjd'
$ws.Range("D3").Value = 'skksjds'
$ws.Range("E3").Formula = '=TEXT(55,"0")'
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4163) | Out-Null

$ws.Range("A4").Value = 'fg'
$ws.Range("B4").Value = 'gsfg'
$ws.Range("C4").Value = 'df'
$ws.Range("D4").Value = 'cvc'
$ws.Range("E4").Formula = '=TEXT(30,"0")'
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4163) | Out-Null

$ws.Range("A5").Value = 'lkj'
$ws.Range("B5").Value = 'ljlj'
$ws.Range("C5").Value = 'ljlkj'
$ws.Range("D5").Value = 'ljlj'
$ws.Range("E5").Formula = '=TEXT(34,"0")'
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4163) | Out-Null

$ws.Range("A6").Value = 'jujhg'
$ws.Range("B6").Value = 'jg'
$ws.Range("C6").Value = 'jg'
$ws.Range("D6").Formula = "'"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Formula = '=TEXT(35,"0")'
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4163) | Out-Null

$ws.Range("A7").Value = 'sfs'
$ws.Range("B7").Value = 'dfg'
$ws.Range("C7").Value = 'dfg'
$ws.Range("D7").Formula = "'"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Formula = '=TEXT(50,"0")'
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4163) | Out-Null

$ws.Range("A8").Value = 'jh'
$ws.Range("B8").Value = 'jhj'
$ws.Range("C8").Value = 'jhjh'
$ws.Range("D8").Value = 'jh'
$ws.Range("E8").Formula = '=TEXT(52,"0")'
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4163) | Out-Null

$ws.Range("A9").Value = 'df'
$ws.Range("B9").Value = 'df'
$ws.Range("C9").Value = 'df'
$ws.Range("D9").Value = 'df'
$ws.Range("E9").Formula = '=TEXT(48,"0")'
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4163) | Out-Null

$ws.Range("A10").Value = 'dfg'
$ws.Range("B10").Value = 'dfg'
$ws.Range("C10").Value = 'dfg'
$ws.Range("D10").Value = 'fdg'
$ws.Range("E10").Formula = '=TEXT(29,"0")'
$ws.Range("E10").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4163) | Out-Null

$ws.Range("A11").Value = ' sd'
$ws.Range("B11").Value = 'af'
$ws.Range("C11").Value = 'af'
$ws.Range("D11").Value = 'asdf'
$ws.Range("E11").Formula = '=TEXT(47,"0")'
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4163) | Out-Null

$ws.Range("A12").Value = 'sdg'
$ws.Range("B12").Value = 'dfg'
$ws.Range("C12").Value = 'dg'
$ws.Range("D12").Formula = "'"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Formula = '=TEXT(12,"0")'
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4163) | Out-Null

$ws.Range("A13").Value = 'xh'
$ws.Range("B13").Value = 'dfh'
$ws.Range("C13").Value = 'fgh'
$ws.Range("D13").Formula = "'"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Formula = '=TEXT(11,"0")'
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false
